$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was inserted at row 3, pushing the existing
# data rows (previously 3-20) down to rows 4-21.
$ws.Rows("3:3").Insert()

# Populate the newly inserted row 3 with this week's reading.
$ws.Range("A3").Value = 11
$ws.Range("B3").Value = "Vega Monumental Concepción"
$ws.Range("C3").Value = "Bíobío"
$ws.Range("D3").Value = 44545
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 100112022
$ws.Range("G3").Value = "Arveja Verde"
$ws.Range("H3").Value = "Perfection"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 180
$ws.Range("K3").Value = 15000
$ws.Range("L3").Value = 16000
$ws.Range("M3").Value = 15444
$ws.Range("N3").Value = "$/saco 25 kilos"
$ws.Range("O3").Value = "Carahue"
$ws.Range("P3").Value = 618
$ws.Range("Q3").Value = 25
$ws.Range("R3").Value = "Hortaliza"
